{"js": "// Replace the Start time, End time, and Duration values in the document body.\n\nconst replacements = [\n  { search: \"Start time: 2017-12-27 18:30:48\", replace: \"Start time: 2018-01-31 12:33:39\" },\n  { search: \"End time: 2017-12-27 18:31:42\", replace: \"End time: 2018-01-31 12:34:23\" },\n  { search: \"Duration: 53.98 secs\", replace: \"Duration: 44.14 secs\" }\n];\n\nfor (const { search, replace } of replacements) {\n  const results = context.document.body.search(search, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the Start time, End time, and Duration values in the document body.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Start time: 2017-12-27 18:30:48\"; Replace = \"Start time: 2018-01-31 12:33:39\" },\n    @{ Find = \"End time: 2017-12-27 18:31:42\"; Replace = \"End time: 2018-01-31 12:34:23\" },\n    @{ Find = \"Duration: 53.98 secs\"; Replace = \"Duration: 44.14 secs\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
